$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '301.55'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '-0.52%'
$c.Style = 'Normal'

# Row 3
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '31.45'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '-1.83%'
$c.Style = 'Normal'

# Row 4
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '5.155'
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '-1.99%'
$c.Style = 'Normal'

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '0.07374'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '-1.16%'
$c.Style = 'Normal'

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '2.362'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '53.93%'
$c.Style = 'Normal'

# Row 7
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '7.918'
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '0.86%'
$c.Style = 'Normal'

# Row 8
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '3.758'
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '-0.85%'
$c.Style = 'Normal'

# Row 9
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.9226'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '0.35%'
$c.Style = 'Normal'

# Row 10
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.1747'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '3.70%'
$c.Style = 'Normal'

# Row 11
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07491'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '-6.40%'
$c.Style = 'Normal'

# Row 12
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.08134'
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '2.25%'
$c.Style = 'Normal'

# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.03038'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '0.24%'
$c.Style = 'Normal'

# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.09944'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '0.25%'
$c.Style = 'Normal'

# Row 15
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.001495'
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '0.23%'
$c.Style = 'Normal'

# Row 16
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.006106'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '-3.41%'
$c.Style = 'Normal'

# Row 17
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '3.452'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '-0.55%'
$c.Style = 'Normal'

# Row 18
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.226'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '-0.36%'
$c.Style = 'Normal'

# Row 19
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '-1.03%'
$c.Style = 'Normal'

# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.1339'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '-0.37%'
$c.Style = 'Normal'

# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '4.652'
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '3.77%'
$c.Style = 'Normal'

# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.04634'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '0.62%'
$c.Style = 'Normal'

# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.1569'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '-3.21%'
$c.Style = 'Normal'

# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.001225'
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '0.48%'
$c.Style = 'Normal'

# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.004484'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '0.85%'
$c.Style = 'Normal'

# Row 26
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '-7.13%'
$c.Style = 'Normal'

# Row 27
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '6.99%'
$c.Style = 'Normal'

# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01725'
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '-1.13%'
$c.Style = 'Normal'

# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.04526'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '0.49%'
$c.Style = 'Normal'

# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.007248'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '0.92%'
$c.Style = 'Normal'

# Row 42
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.1347'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '-0.24%'
$c.Style = 'Normal'

# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.002221'
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '1.38%'
$c.Style = 'Normal'

# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.01091'
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '-15.02%'
$c.Style = 'Normal'

# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.00006293'
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '2.57%'
$c.Style = 'Normal'

# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.01001'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '-23.05%'
$c.Style = 'Normal'

# Row 47
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.7106'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '-61.96%'
$c.Style = 'Normal'
